$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap week 5 / week 6 topics (midterm moved from week 5 to week 6) ---
$c5Value = $ws.Range("C5").Value2
$c6Value = $ws.Range("C6").Value2
$c5Wrap  = $ws.Range("C5").WrapText
$c6Wrap  = $ws.Range("C6").WrapText

$ws.Range("C5").Value2 = $c6Value
$ws.Range("C6").Value2 = $c5Value

if ($c6Wrap) { $ws.Range("C5").WrapText = $true } else { $ws.Range("C5").Style = "Normal" }
if ($c5Wrap) { $ws.Range("C6").WrapText = $true } else { $ws.Range("C6").Style = "Normal" }

# row 5 had the taller (17pt) formatting, row 6 had the default height;
# the swap moves the taller formatting onto row 6 and row 5 reverts to default
$ws.Rows(5).AutoFit()
$ws.Rows(6).RowHeight = 17

# --- Swap week 11 / week 12 topics (midterm moved from week 11 to week 12) ---
$c11Value = $ws.Range("C11").Value2
$c12Value = $ws.Range("C12").Value2
$c11Wrap  = $ws.Range("C11").WrapText
$c12Wrap  = $ws.Range("C12").WrapText

$ws.Range("C11").Value2 = $c12Value
$ws.Range("C12").Value2 = $c11Value

if ($c12Wrap) { $ws.Range("C11").WrapText = $true } else { $ws.Range("C11").Style = "Normal" }
if ($c11Wrap) { $ws.Range("C12").WrapText = $true } else { $ws.Range("C12").Style = "Normal" }

# row 11 had the taller (17pt) formatting, row 12 had the default height;
# the swap moves the taller formatting onto row 12 and row 11 reverts to default
$ws.Rows(11).AutoFit()
$ws.Rows(12).RowHeight = 17

# --- Update the active selection to reflect where the user ended up ---
$ws.Range("C11").Select()
